# Phase Two Time Log update:
#  - Remove the old "Phase 2 Requirements" log entry (row 9) and fold its
#    start-time data up into row 2, now logged against the new
#    "System Architecture Document" activity.
#  - Add a new "UML Diagram" log entry in row 3.
#  - Clear out the now-unused date placeholders that used to occupy rows 3-9
#    (and their shared D-column formulas).
#  - Trim the trailing empty rows at the bottom of the sheet (22-28) and tidy
#    up row 21.
#  - Move the active selection to F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Break up the D2:D8 shared formula group and drop the formulas that are no
# longer needed (rows 4-9); this also removes those cells entirely, matching
# a fully blank row.
$ws.Range("D4:D9").Clear()

# Clear the leftover date values from rows 3-9 (keep the cell styling).
$ws.Range("A3:A9").ClearContents()
$ws.Range("B9:C9").ClearContents()
$ws.Range("E9").ClearContents()

# Row 2: "System Architecture Document" work, using what used to be row 9's
# start time, plus a new end time.
$ws.Range("A2").Value = 44490
$ws.Range("B2").Value = 0.74305555555555547
$ws.Range("C2").Value = 0.75
$ws.Range("D2").Formula = "=C2-B2"
$ws.Range("E2").Value = "System Architecture Document"
$ws.Range("F2").Value = "Started System Architecture Document."

# Row 3: new "UML Diagram" entry.
$ws.Range("B3").Value = 0.75
$ws.Range("C3").Value = 0.83333333333333337
$ws.Range("D3").Formula = "=C3-B3"
$ws.Range("E3").Value = "UML Diagram"
$ws.Range("F3").Value = "Creating UML class diagrams for architecture document."

# Tidy row 21 (drop A21/C21, keep B21) and drop the now-unused trailing rows.
$ws.Range("A21").Clear()
$ws.Range("C21").Clear()
$ws.Range("A22:C28").Clear()

# Match the saved selection from the edit.
$ws.Range("F3").Select()
